$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Rent -> ferrari
$ws.Range("A2").Value = "ferrari"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = 46325.041666666664

# Row 3: Food stays, amount and date change
$ws.Range("B3").Value = 120
$ws.Range("C3").Value = 46138.083333333336

# Row 4: Stock: Ethereum -> But business with bestfriendo
$ws.Range("A4").Value = "But business with bestfriendo"
$ws.Range("B4").Value = 1300
$ws.Range("C4").Value = 45935.083333333336

# Row 5: Fuel -> food
$ws.Range("A5").Value = "food"
$ws.Range("B5").Value = 300
$ws.Range("C5").Value = 45930.083333333336
